$d = $word.ActiveDocument

$oldFull = "照片:红包给最亲爱的你"
$count = $d.Paragraphs.Count

for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $r = $p.Range
    $text = $r.Text

    if ($text.Length -ge $oldFull.Length -and $text.Substring(0, $oldFull.Length) -eq $oldFull) {
        $base = $r.Start

        # Remove "红包" (characters at offsets 3-4, i.e. positions base+3..base+5)
        $redPacket = $d.Range($base + 3, $base + 5)
        $redPacket.Text = ""

        # Remove "亲" from "给最亲爱的你" (now at offset 5, i.e. position base+5..base+6)
        $qin = $d.Range($base + 5, $base + 6)
        $qin.Text = ""
    }
}
